$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "last updated" timestamp banner (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 20 de Mayo de 2020 a las 10:05"

# --- Refresh country stats (Rusia, row 5) ---
$ws.Range("B5").Value = 308705
$ws.Range("C5").Value = 8764
$ws.Range("D5").Value = 85392
$ws.Range("E5").Value = 220341
$ws.Range("G5").Value = 135
$ws.Range("H5").Value = 2972

# --- Refresh country stats (Ucrania, row 35) ---
$ws.Range("B35").Value = 19230
$ws.Range("C35").Value = 354
$ws.Range("D35").Value = 5955
$ws.Range("E35").Value = 12711
$ws.Range("G35").Value = 16
$ws.Range("H35").Value = 564

# --- Refresh country stats (Eslovaquia, row 97) ---
$ws.Range("B97").Value = 1496
$ws.Range("C97").Value = 1
$ws.Range("D97").Value = 1231
$ws.Range("E97").Value = 237

# --- Reorder Belice / Nueva Caledonia / Santa Lucia (rows 195-197) ---
# Previously: 195=Belice, 196=Nueva Caledonia, 197=Santa Lucia
# Now:        195=Nueva Caledonia, 196=Santa Lucia, 197=Belice
# (each country keeps its own stats, only the row order changes)
$ws.Range("A195").Value = "Nueva Caledonia"
$ws.Range("B195").Value = 18
$ws.Range("C195").Value = 0
$ws.Range("D195").Value = 18
$ws.Range("E195").Value = 0
$ws.Range("F195").Value = 0
$ws.Range("G195").Value = 0
$ws.Range("H195").Value = 0

$ws.Range("A196").Value = "Santa Lucia"
$ws.Range("B196").Value = 18
$ws.Range("C196").Value = 0
$ws.Range("D196").Value = 18
$ws.Range("E196").Value = 0
$ws.Range("F196").Value = 0
$ws.Range("G196").Value = 0
$ws.Range("H196").Value = 0

$ws.Range("A197").Value = "Belice"
$ws.Range("B197").Value = 18
$ws.Range("C197").Value = 0
$ws.Range("D197").Value = 16
$ws.Range("E197").Value = 0
$ws.Range("F197").Value = 0
$ws.Range("G197").Value = 0
$ws.Range("H197").Value = 2
